# Apply updated results (Results Updated and PathEnrResults Added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 10,18
$arr[0,0] = 253.96551724137922
$arr[0,1] = 0.9693486590038315
$arr[0,2] = 0.9616858237547893
$arr[0,3] = 0.9731800766283526
$arr[0,4] = 0.9547456375042581
$arr[0,5] = 0.9951309067688379
$arr[0,6] = 0.9616858237547893
$arr[0,7] = 0.955322933771209
$arr[0,8] = 0.9317821508826001
$arr[0,9] = 84.32515256584249
$arr[0,10] = 0.04129560937376906
$arr[0,11] = 0.07916370629017346
$arr[0,12] = 0.05002301291071051
$arr[0,13] = 0.07916370629017346
$arr[0,14] = 0.08047560872662458
$arr[0,15] = 0.05987099927422121
$arr[0,16] = 0.09071061184605762
$arr[0,17] = 0.012669100181615926
$arr[1,0] = 239.89108910891093
$arr[1,1] = 0.9735973597359733
$arr[1,2] = 0.9653465346534653
$arr[1,3] = 0.9777227722772274
$arr[1,4] = 0.9610600620501611
$arr[1,5] = 0.9942244224422438
$arr[1,6] = 0.9653465346534653
$arr[1,7] = 0.9629302215935874
$arr[1,8] = 0.9412521787295782
$arr[1,9] = 76.82628469346922
$arr[1,10] = 0.03896346769966366
$arr[1,11] = 0.07194454532318967
$arr[1,12] = 0.046979737223297936
$arr[1,13] = 0.07194454532318967
$arr[1,14] = 0.07559233772144715
$arr[1,15] = 0.055872217072821
$arr[1,16] = 0.08521684068884133
$arr[1,17] = 0.01546817306039549
$arr[2,0] = 228.11650485436905
$arr[2,1] = 0.9746494066882411
$arr[2,2] = 0.9660194174757281
$arr[2,3] = 0.9789644012944982
$arr[2,4] = 0.9624563139126245
$arr[2,5] = 0.9948759439050695
$arr[2,6] = 0.9660194174757281
$arr[2,7] = 0.9648636153490522
$arr[2,8] = 0.9434683425481231
$arr[2,9] = 78.59979306599435
$arr[2,10] = 0.03712682392526563
$arr[2,11] = 0.07139759223530995
$arr[2,12] = 0.04476768677440194
$arr[2,13] = 0.07139759223530995
$arr[2,14] = 0.07314652532981213
$arr[2,15] = 0.05426269315081091
$arr[2,16] = 0.08204729164060874
$arr[2,17] = 0.014649681944059462
$arr[3,0] = 211.13761467889918
$arr[3,1] = 0.9765545361875636
$arr[3,2] = 0.9663608562691128
$arr[3,3] = 0.9816513761467887
$arr[3,4] = 0.9650102191386597
$arr[3,5] = 0.9949031600407747
$arr[3,6] = 0.9663608562691128
$arr[3,7] = 0.9688728702490171
$arr[3,8] = 0.9475046652589532
$arr[3,9] = 71.47823650774801
$arr[3,10] = 0.03648788718412104
$arr[3,11] = 0.07092791983137557
$arr[3,12] = 0.04144663844159234
$arr[3,13] = 0.07092791983137557
$arr[3,14] = 0.06886394290698962
$arr[3,15] = 0.05387912599358178
$arr[3,16] = 0.08109649833982617
$arr[3,17] = 0.014444957658776687
$arr[4,0] = 200.24000000000018
$arr[4,1] = 0.9755555555555552
$arr[4,2] = 0.9649999999999999
$arr[4,3] = 0.9808333333333327
$arr[4,4] = 0.9637212787212788
$arr[4,5] = 0.995486111111111
$arr[4,6] = 0.9649999999999999
$arr[4,7] = 0.9674999999999999
$arr[4,8] = 0.9454030100334446
$arr[4,9] = 72.32797354810765
$arr[4,10] = 0.03646015374110923
$arr[4,11] = 0.06822669679005593
$arr[4,12] = 0.0424670599223677
$arr[4,13] = 0.06822669679005593
$arr[4,14] = 0.07047737265495682
$arr[4,15] = 0.05296267831786668
$arr[4,16] = 0.0804535837704723
$arr[4,17] = 0.014076143687523917
$arr[5,0] = 180.63999999999996
$arr[5,1] = 0.9766666666666667
$arr[5,2] = 0.9666666666666665
$arr[5,3] = 0.9816666666666665
$arr[5,4] = 0.9650729270729269
$arr[5,5] = 0.9956249999999995
$arr[5,6] = 0.9666666666666665
$arr[5,7] = 0.9692857142857142
$arr[5,8] = 0.9476834296138642
$arr[5,9] = 69.69913988716146
$arr[5,10] = 0.033668350126272215
$arr[5,11] = 0.07106690545187216
$arr[5,12] = 0.04030606808631133
$arr[5,13] = 0.07106690545187216
$arr[5,14] = 0.06657770811845104
$arr[5,15] = 0.04997485691631375
$arr[5,16] = 0.07496879081790679
$arr[5,17] = 0.013485117112067268
$arr[6,0] = 158.82000000000002
$arr[6,1] = 0.9761111111111106
$arr[6,2] = 0.9649999999999997
$arr[6,3] = 0.9816666666666665
$arr[6,4] = 0.9640772560772557
$arr[6,5] = 0.9952777777777773
$arr[6,6] = 0.9649999999999997
$arr[6,7] = 0.9684523809523807
$arr[6,8] = 0.9462954697476436
$arr[6,9] = 73.35817485585241
$arr[6,10] = 0.03464947485885504
$arr[6,11] = 0.07222222222222482
$arr[6,12] = 0.03852644958492143
$arr[6,13] = 0.07222222222222482
$arr[6,14] = 0.06565086674282744
$arr[6,15] = 0.05194011360887088
$arr[6,16] = 0.07761785866715273
$arr[6,17] = 0.012977985523274735
$arr[7,0] = 128.10000000000002
$arr[7,1] = 0.9733333333333328
$arr[7,2] = 0.9616666666666667
$arr[7,3] = 0.9791666666666661
$arr[7,4] = 0.9602707292707291
$arr[7,5] = 0.9952083333333329
$arr[7,6] = 0.9616666666666667
$arr[7,7] = 0.964880952380952
$arr[7,8] = 0.9403577885882232
$arr[7,9] = 65.89531029740978
$arr[7,10] = 0.03826255278343741
$arr[7,11] = 0.07436600722307987
$arr[7,12] = 0.044907624352969515
$arr[7,13] = 0.07436600722307987
$arr[7,14] = 0.07256935351705382
$arr[7,15] = 0.055814239468645326
$arr[7,16] = 0.08440829874482399
$arr[7,17] = 0.013558607996074425
$arr[8,0] = 90.26999999999997
$arr[8,1] = 0.9744444444444442
$arr[8,2] = 0.9616666666666668
$arr[8,3] = 0.9808333333333329
$arr[8,4] = 0.9613616383616387
$arr[8,5] = 0.993958333333333
$arr[8,6] = 0.9616666666666668
$arr[8,7] = 0.9667857142857137
$arr[8,8] = 0.9423867740954698
$arr[8,9] = 59.06867343878802
$arr[8,10] = 0.03295709855904684
$arr[8,11] = 0.07436600722307891
$arr[8,12] = 0.0371830036115245
$arr[8,13] = 0.07436600722307891
$arr[8,14] = 0.06379112019502256
$arr[8,15] = 0.050223581316696686
$arr[8,16] = 0.07444368819542706
$arr[8,17] = 0.015375510022579296
$arr[9,0] = 45.12000000000001
$arr[9,1] = 0.9649999999999999
$arr[9,2] = 0.9516666666666663
$arr[9,3] = 0.9716666666666665
$arr[9,4] = 0.9484775224775224
$arr[9,5] = 0.9897222222222218
$arr[9,6] = 0.9516666666666663
$arr[9,7] = 0.9529166666666663
$arr[9,8] = 0.922206283571501
$arr[9,9] = 39.832802075828425
$arr[9,10] = 0.0444269045075346
$arr[9,11] = 0.0796139676086065
$arr[9,12] = 0.054561819467998825
$arr[9,13] = 0.0796139676086065
$arr[9,14] = 0.08601344948357957
$arr[9,15] = 0.06297385193042307
$arr[9,16] = 0.09656761895675348
$arr[9,17] = 0.024685482606220884

$ws.Range("B2:S11").Value = $arr

